$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
# D5 (MyForecast for week W4): 10 -> 11
$wsForecast.Range("D5").Value = 11

# --- Sheet: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")
# These cells store numeric-looking values as text, so force text format
# before assigning to keep them as strings rather than numbers.
$wsSummary.Range("B9:B11").NumberFormat = "@"
# B9 (Total Forecast 16 Weeks): "238" -> "239"
$wsSummary.Range("B9").Value = "239"
# B10 (Total Forecast 8 Weeks): "108" -> "109"
$wsSummary.Range("B10").Value = "109"
# B11 (Total Forecast 4 Weeks): "48" -> "49"
$wsSummary.Range("B11").Value = "49"
